$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values from the refreshed symbol list.
# Values are written with a leading apostrophe so Excel keeps them as literal text
# (matching the original inlineStr cells) instead of coercing to numbers/percentages.
$ws.Range("D2").Value = "'280.97"
$ws.Range("E2").Value = "'5.76%"
$ws.Range("E3").Value = "'0.20%"
$ws.Range("D4").Value = "'4.943"
$ws.Range("E4").Value = "'5.22%"
$ws.Range("D5").Value = "'0.06414"
$ws.Range("E5").Value = "'5.43%"
$ws.Range("D6").Value = "'6.987"
$ws.Range("E6").Value = "'3.78%"
$ws.Range("E7").Value = "'5.85%"
$ws.Range("D8").Value = "'0.8870"
$ws.Range("E8").Value = "'4.27%"
$ws.Range("D9").Value = "'1.016"
$ws.Range("E9").Value = "'11.72%"
$ws.Range("D10").Value = "'0.1489"
$ws.Range("E10").Value = "'5.84%"
$ws.Range("D11").Value = "'0.05207"
$ws.Range("E11").Value = "'3.23%"
$ws.Range("D12").Value = "'0.07425"
$ws.Range("E12").Value = "'4.62%"
$ws.Range("D13").Value = "'0.03110"
$ws.Range("E13").Value = "'-1.27%"
$ws.Range("D14").Value = "'0.09046"
$ws.Range("E14").Value = "'0.23%"
$ws.Range("D15").Value = "'0.001564"
$ws.Range("E15").Value = "'2.34%"
$ws.Range("D16").Value = "'0.0006333"
$ws.Range("E16").Value = "'3.94%"
$ws.Range("D17").Value = "'0.006009"
$ws.Range("E17").Value = "'0.61%"
$ws.Range("D18").Value = "'3.495"
$ws.Range("E18").Value = "'1.12%"
$ws.Range("E19").Value = "'5.66%"
$ws.Range("E21").Value = "'3.82%"
$ws.Range("D22").Value = "'3.933"
$ws.Range("E22").Value = "'-3.75%"
$ws.Range("D23").Value = "'0.04353"
$ws.Range("E23").Value = "'2.54%"
$ws.Range("D24").Value = "'0.001180"
$ws.Range("E24").Value = "'0.09%"
$ws.Range("D25").Value = "'0.003688"
$ws.Range("E25").Value = "'-9.13%"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("D27").Value = "'0.0001694"
$ws.Range("E27").Value = "'0.72%"
$ws.Range("E40").Value = "'4.43%"
$ws.Range("D41").Value = "'0.006655"
$ws.Range("E41").Value = "'58.78%"
$ws.Range("D42").Value = "'0.1177"
$ws.Range("E42").Value = "'5.60%"
$ws.Range("E43").Value = "'11.78%"
$ws.Range("E44").Value = "'12.64%"
$ws.Range("D45").Value = "'0.00005242"
$ws.Range("E45").Value = "'2.57%"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("E47").Value = "'813.01%"
$ws.Range("D48").Value = "'0.02250"
$ws.Range("E48").Value = "'-8.06%"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("E50").Value = "'-0.17%"
